$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.255269829919371
$ws.Range("C2").Value = 0.1762097630474386
$ws.Range("D2").Value = 0.1273133655189014
$ws.Range("E2").Value = 0.1219173530156605
$ws.Range("F2").Value = 1.582352370776682
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.9683392973090825
$ws.Range("J2").Value = 0.1535881532749244
$ws.Range("L2").Value = 0.3136387744265221
$ws.Range("M2").Value = 0.3075114984643221
$ws.Range("O2").Value = 4.051546183316731
$ws.Range("B3").Value = 1.151776496250079
$ws.Range("C3").Value = 0.1600266875665852
$ws.Range("D3").Value = 0.1269731100004137
$ws.Range("E3").Value = 0.1227890251941952
$ws.Range("F3").Value = 1.594128176722563
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.9813925671496371
$ws.Range("J3").Value = 0.1545577064332466
$ws.Range("L3").Value = 0.3096592936175568
$ws.Range("M3").Value = 0.2905356750175656
$ws.Range("O3").Value = 4.086702161348327
$ws.Range("B4").Value = 1.088267293464526
$ws.Range("C4").Value = 0.1500344830694473
$ws.Range("D4").Value = 0.1267937272371427
$ws.Range("E4").Value = 0.1233555329943212
$ws.Range("F4").Value = 1.602255175343331
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.9899683364320921
$ws.Range("J4").Value = 0.1551909129100846
$ws.Range("L4").Value = 0.307311879447667
$ws.Range("M4").Value = 0.2801605470758588
$ws.Range("O4").Value = 4.110764405732596
$ws.Range("B5").Value = 1.062397676821718
$ws.Range("C5").Value = 0.1459488057901694
$ws.Range("D5").Value = 0.1267280878345645
$ws.Range("E5").Value = 0.1235942730177677
$ws.Range("F5").Value = 1.605792459628823
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.993604038086449
$ws.Range("J5").Value = 0.1554584992766248
$ws.Range("L5").Value = 0.3063795564471761
$ws.Range("M5").Value = 0.2759450119319098
$ws.Range("O5").Value = 4.121192419656083
$ws.Range("B6").Value = 1.058102758662926
$ws.Range("C6").Value = 0.1452695574827203
$ws.Range("D6").Value = 0.1267176400486321
$ws.Range("E6").Value = 0.1236343921564678
$ws.Range("F6").Value = 1.606393440824036
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.9942162581111873
$ws.Range("J6").Value = 0.1555035091536894
$ws.Range("L6").Value = 0.3062262143442069
$ws.Range("M6").Value = 0.2752457847349135
$ws.Range("O6").Value = 4.122961571078477
$ws.Range("B7").Value = 1.087918360221806
$ws.Range("C7").Value = 0.1499794376406953
$ws.Range("D7").Value = 0.1267928117479755
$ws.Range("E7").Value = 0.1233587207878397
$ws.Range("F7").Value = 1.602301967483548
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.9900167979344694
$ws.Range("J7").Value = 0.1551944829810772
$ws.Range("L7").Value = 0.3072992073874588
$ws.Range("M7").Value = 0.2801036442068821
$ws.Range("O7").Value = 4.1109025216075
$ws.Range("B8").Value = 1.219578884434895
$ws.Range("C8").Value = 0.1706415448983876
$ws.Range("D8").Value = 0.1271899351616241
$ws.Range("E8").Value = 0.122211417934631
$ws.Range("F8").Value = 1.586226581612159
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.9727235898947804
$ws.Range("J8").Value = 0.1539146017823541
$ws.Range("L8").Value = 0.3122467980090917
$ws.Range("M8").Value = 0.3016484377177022
$ws.Range("O8").Value = 4.06315385639715
$ws.Range("B9").Value = 1.477984808978306
$ws.Range("C9").Value = 0.2107090024373974
$ws.Range("D9").Value = 0.1282016814555931
$ws.Range("E9").Value = 0.1202092679808198
$ws.Range("F9").Value = 1.561817346997117
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.943265610404989
$ws.Range("J9").Value = 0.1517045478677037
$ws.Range("L9").Value = 0.3227059155219223
$ws.Range("M9").Value = 0.3442675539320419
$ws.Range("O9").Value = 3.989179689005738
$ws.Range("B10").Value = 1.667899729087026
$ws.Range("C10").Value = 0.2398628693825344
$ws.Range("D10").Value = 0.129085382618328
$ws.Range("E10").Value = 0.1188884091111198
$ws.Range("F10").Value = 1.548222672724989
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.924341245825989
$ws.Range("J10").Value = 0.1502623289237306
$ws.Range("L10").Value = 0.3308462316967962
$ws.Range("M10").Value = 0.3757927740185778
$ws.Range("O10").Value = 3.94683564734774
$ws.Range("B11").Value = 1.754297375352394
$ws.Range("C11").Value = 0.253062443140692
$ws.Range("D11").Value = 0.1295175412505074
$ws.Range("E11").Value = 0.1183199238025776
$ws.Range("F11").Value = 1.54298065069338
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.9163230144484018
$ws.Range("J11").Value = 0.1496453755646154
$ws.Range("L11").Value = 0.3346474127728385
$ws.Range("M11").Value = 0.3901782338862034
$ws.Range("O11").Value = 3.930183228263871
$ws.Range("B12").Value = 1.787013009736484
$ws.Range("C12").Value = 0.2580515661375102
$ws.Range("D12").Value = 0.1296854940011229
$ws.Range("E12").Value = 0.1181092955184613
$ws.Range("F12").Value = 1.541131153565296
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.9133717140742164
$ws.Range("J12").Value = 0.1494173569158797
$ws.Range("L12").Value = 0.3361008179275728
$ws.Range("M12").Value = 0.3956317567994887
$ws.Range("O12").Value = 3.924253051532332
$ws.Range("B13").Value = 1.779967199694397
$ws.Range("C13").Value = 0.2569774855626861
$ws.Range("D13").Value = 0.1296491314901473
$ws.Range("E13").Value = 0.1181544516924378
$ws.Range("F13").Value = 1.541523446828116
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.9140035459651727
$ws.Range("J13").Value = 0.1494662156557229
$ws.Range("L13").Value = 0.3357871813700513
$ws.Range("M13").Value = 0.3944569787830901
$ws.Range("O13").Value = 3.925513503966499
$ws.Range("B14").Value = 1.756988946819206
$ws.Range("C14").Value = 0.2534730885275849
$ws.Range("D14").Value = 0.1295312727623568
$ws.Range("E14").Value = 0.1183025022641555
$ws.Range("F14").Value = 1.54282577470795
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.9160785044452986
$ws.Range("J14").Value = 0.1496265040338933
$ws.Range("L14").Value = 0.3347667058315409
$ws.Range("M14").Value = 0.390626778926368
$ws.Range("O14").Value = 3.929687816502707
$ws.Range("B15").Value = 1.74291388036886
$ws.Range("C15").Value = 0.2513253294292781
$ws.Range("D15").Value = 0.1294596403706478
$ws.Range("E15").Value = 0.1183937920613474
$ws.Range("F15").Value = 1.543641141182775
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.9173605530378381
$ws.Range("J15").Value = 0.1497254151431182
$ws.Range("L15").Value = 0.3341434519742847
$ws.Range("M15").Value = 0.3882814504083569
$ws.Range("O15").Value = 3.932293646929679
$ws.Range("B16").Value = 1.662253364477863
$ws.Range("C16").Value = 0.2389989631051037
$ws.Range("D16").Value = 0.1290577440683833
$ws.Range("E16").Value = 0.1189262115518979
$ws.Range("F16").Value = 1.548584215939428
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.9248771488939944
$ws.Range("J16").Value = 0.1503034339891105
$ws.Range("L16").Value = 0.3305997788699244
$ws.Range("M16").Value = 0.3748535167683471
$ws.Range("O16").Value = 3.947976478545286
$ws.Range("B17").Value = 1.612770512264206
$ws.Range("C17").Value = 0.2314209021710667
$ws.Range("D17").Value = 0.1288188944185222
$ws.Range("E17").Value = 0.1192611187981054
$ws.Range("F17").Value = 1.551858001861191
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.9296396700124703
$ws.Range("J17").Value = 0.1506680369693161
$ws.Range("L17").Value = 0.3284508875233882
$ws.Range("M17").Value = 0.3666270678350898
$ws.Range("O17").Value = 3.958266184863817
$ws.Range("B18").Value = 1.584309811644175
$ws.Range("C18").Value = 0.2270563187604751
$ws.Range("D18").Value = 0.1286843550517887
$ws.Range("E18").Value = 0.1194567970282021
$ws.Range("F18").Value = 1.553829684505814
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.9324345257347382
$ws.Range("J18").Value = 0.1508814300583179
$ws.Range("L18").Value = 0.3272241442763431
$ws.Range("M18").Value = 0.3618996384379329
$ws.Range("O18").Value = 3.964430202551739
$ws.Range("B19").Value = 1.574673656450614
$ws.Range("C19").Value = 0.225577543252399
$ws.Range("D19").Value = 0.1286392912262997
$ws.Range("E19").Value = 0.1195235742369887
$ws.Range("F19").Value = 1.554512493118089
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.9333903583950729
$ws.Range("J19").Value = 0.1509543143930134
$ws.Range("L19").Value = 0.3268103813541927
$ws.Range("M19").Value = 0.3602997447643119
$ws.Range("O19").Value = 3.966559411436748
$ws.Range("B20").Value = 1.618038005682308
$ws.Range("C20").Value = 0.2322282104556734
$ws.Range("D20").Value = 0.1288440266177986
$ws.Range("E20").Value = 0.1192251519310306
$ws.Range("F20").Value = 1.551500322588701
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.9291269384477872
$ws.Range("J20").Value = 0.1506288432969942
$ws.Range("L20").Value = 0.3286786850686383
$ws.Range("M20").Value = 0.3675023544466995
$ws.Range("O20").Value = 3.957145401185358
$ws.Range("B21").Value = 1.763738266572034
$ws.Range("C21").Value = 0.2545026684908009
$ws.Range("D21").Value = 0.1295657741898921
$ws.Range("E21").Value = 0.1182588902556874
$ws.Range("F21").Value = 1.542439570480575
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.9154667302930548
$ws.Range("J21").Value = 0.149579271382148
$ws.Range("L21").Value = 0.3350660658258988
$ws.Range("M21").Value = 0.391751639618974
$ws.Range("O21").Value = 3.928451519232539
$ws.Range("B22").Value = 1.858953476297529
$ws.Range("C22").Value = 0.2690061687768832
$ws.Range("D22").Value = 0.1300625477050943
$ws.Range("E22").Value = 0.1176544503777204
$ws.Range("F22").Value = 1.537307913591277
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.9070346637577096
$ws.Range("J22").Value = 0.1489259965597807
$ws.Range("L22").Value = 0.3393220178543856
$ws.Range("M22").Value = 0.4076351264063831
$ws.Range("O22").Value = 3.911888653491161
$ws.Range("B23").Value = 1.808136717429932
$ws.Range("C23").Value = 0.2612704166854485
$ws.Range("D23").Value = 0.1297951273403157
$ws.Range("E23").Value = 0.1179745783550907
$ws.Range("F23").Value = 1.539974469008953
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.9114896267981933
$ws.Range("J23").Value = 0.1492716768964808
$ws.Range("L23").Value = 0.3370431274244652
$ws.Range("M23").Value = 0.3991547018611357
$ws.Range("O23").Value = 3.920528026319602
$ws.Range("B24").Value = 1.615656610231724
$ws.Range("C24").Value = 0.2318632508308269
$ws.Range("D24").Value = 0.1288326556939552
$ws.Range("E24").Value = 0.1192414027838256
$ws.Range("F24").Value = 1.55166175047998
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.9293585674072489
$ws.Range("J24").Value = 0.1506465509865365
$ws.Range("L24").Value = 0.3285756707429499
$ws.Range("M24").Value = 0.36710663091025
$ws.Range("O24").Value = 3.957651333982994
$ws.Range("B25").Value = 1.408062899650076
$ws.Range("C25").Value = 0.199918883432133
$ws.Range("D25").Value = 0.1279032070290995
$ws.Range("E25").Value = 0.1207244759200918
$ws.Range("F25").Value = 1.567658865321278
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.9507576211713591
$ws.Range("J25").Value = 0.1522704618070776
$ws.Range("L25").Value = 0.3197959785233593
$ws.Range("M25").Value = 0.3227059155219223
$ws.Range("O25").Value = 4.007084861261006
